$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 801
$ws.Range("F6").Value = 1133
$ws.Range("F8").Value = 39
$ws.Range("F10").Value = 115
$ws.Range("F12").Value = 52
$ws.Range("F14").Value = 798
$ws.Range("F15").Value = 822
$ws.Range("F16").Value = 188
$ws.Range("F17").Value = 49
$ws.Range("F18").Value = 67
$ws.Range("F20").Value = 195
$ws.Range("F22").Value = 2382
$ws.Range("F23").Value = 664
$ws.Range("F25").Value = 1921
$ws.Range("F26").Value = 346
$ws.Range("F27").Value = 2793
$ws.Range("F28").Value = 515
$ws.Range("F30").Value = 685
$ws.Range("F34").Value = 966
$ws.Range("F35").Value = 1698
$ws.Range("F36").Value = 338
$ws.Range("F38").Value = 537
$ws.Range("F39").Value = 158
$ws.Range("F40").Value = 117

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 2
$ws.Range("F12").Value = 72

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 801
$ws.Range("F7").Value = 1133
$ws.Range("F9").Value = 39
$ws.Range("F11").Value = 115
$ws.Range("F13").Value = 52
$ws.Range("F14").Value = 798
$ws.Range("F15").Value = 822
$ws.Range("F16").Value = 188
$ws.Range("F19").Value = 49
$ws.Range("F21").Value = 67
$ws.Range("F22").Value = 195
$ws.Range("F24").Value = 2382
$ws.Range("F25").Value = 664
$ws.Range("F29").Value = 2793
$ws.Range("F30").Value = 515
$ws.Range("F31").Value = 2
$ws.Range("F36").Value = 72
$ws.Range("F37").Value = 685
$ws.Range("F41").Value = 966
$ws.Range("F42").Value = 1698
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "2024-07-13"
$ws.Range("B43").NumberFormat = "General"
$ws.Range("C43").Value = "杭州·代号鸢only-广陵大学"
$ws.Range("D43").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E43").Value = "2024.07.13 09:00-07.13 18:00"
$ws.Range("F43").Value = 338
$ws.Range("G43").Value = 68
$ws.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=83289"
$ws.Range("I43").Value = "//i0.hdslb.com/bfs/openplatform/202403/I3yffJ7Q1711344958258.png"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "2024-07-20"
$ws.Range("B44").NumberFormat = "General"
$ws.Range("C44").Value = "杭州·次元幻想--二次元全女夜场"
$ws.Range("D44").Value = "保淑路2号 The Queen皇后"
$ws.Range("E44").Value = "2024.07.20 13:00-07.20 19:00"
$ws.Range("F44").Value = 537
$ws.Range("G44").Value = 158
$ws.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=81808"
$ws.Range("I44").Value = "//i0.hdslb.com/bfs/openplatform/202402/sUUtSPh91707295826425.jpeg"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "2024-07-27"
$ws.Range("B45").NumberFormat = "General"
$ws.Range("C45").Value = "杭州·夏之誓国乙only-日夜场"
$ws.Range("D45").Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$ws.Range("E45").Value = "2024.07.27 10:00-07.27 21:00"
$ws.Range("F45").Value = 158
$ws.Range("G45").Value = 69
$ws.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=83589"
$ws.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202403/G8PdP81U1711604984731.png"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "2024-08-03"
$ws.Range("B46").NumberFormat = "General"
$ws.Range("C46").Value = "杭州·梦漫星河动漫展"
$ws.Range("D46").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E46").Value = "2024.08.03 10:00-08.04 17:00"
$ws.Range("F46").Value = 117
$ws.Range("G46").Value = 68
$ws.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=82836"
$ws.Range("I46").Value = "//i0.hdslb.com/bfs/openplatform/202403/VFfQUJdD1711700169290.jpeg"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "2024-08-14"
$ws.Range("B47").NumberFormat = "General"
$ws.Range("C47").Value = "杭州·第五幼儿园·第五人格only展"
$ws.Range("D47").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("E47").Value = "2024.08.14 09:00-08.14 17:00"
$ws.Range("F47").Value = 157
$ws.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=82834"
$ws.Range("I47").Value = "//i2.hdslb.com/bfs/openplatform/202403/ftH5TCpR1711444351628.jpeg"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "2024-09-15"
$ws.Range("B48").NumberFormat = "General"
$ws.Range("C48").Value = "杭州·理想乡动漫展-同人创作者大会"
$ws.Range("D48").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws.Range("E48").Value = "2024.09.15 10:00-09.16 17:00"
$ws.Range("F48").Value = 15
$ws.Range("G48").Value = 39
$ws.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=83822"
$ws.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202404/GGEZUjGw1711959030111.png"
